$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2: Global
$t.Cell(2, 2).Range.Text = "157'973'438"
$t.Cell(2, 3).Range.Text = "75'675"
$t.Cell(2, 4).Range.Text = "3'288'455"
$t.Cell(2, 5).Range.Text = "1'373"

# Row 3: USA
$t.Cell(3, 2).Range.Text = "33'515'308"
$t.Cell(3, 3).Range.Text = "35'960"
$t.Cell(3, 4).Range.Text = "596'179"
$t.Cell(3, 5).Range.Text = "334"

# Row 4: Europa
$t.Cell(4, 2).Range.Text = "51'889'020"
$t.Cell(4, 3).Range.Text = "41'967"
$t.Cell(4, 4).Range.Text = "1'097'713"
$t.Cell(4, 5).Range.Text = "1'058"

# Row 6: Frankreich
$t.Cell(6, 2).Range.Text = "5'780'379"
$t.Cell(6, 3).Range.Text = "3'292"
$t.Cell(6, 4).Range.Text = "106'684"
$t.Cell(6, 5).Range.Text = "292"

# Row 7: Deutschland
$t.Cell(7, 2).Range.Text = "3'535'354"
$t.Cell(7, 3).Range.Text = "3'773"
$t.Cell(7, 4).Range.Text = "85'481"
$t.Cell(7, 5).Range.Text = "68"

# Row 9: Schweiz
$t.Cell(9, 4).Range.Text = "10'715"
$t.Cell(9, 5).Range.Text = "5"
